# Refresh the "cryptos" price/volume table with the day's updated figures,
# including a couple of rank swaps (rows 14/15 and 46/47).
#
# NOTE: the D (Price) column cells are plain text in the source file (values
# like "43.734.75" or "0.611" are not valid numbers). Excel's Range.Value
# setter auto-detects numeric-looking strings and would silently convert
# them to numbers (dropping formatting such as trailing zeros / multiple
# "thousand separator" dots). Prefixing with a leading apostrophe forces
# Excel to keep them as text, matching the original cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''43.734.75'
$ws.Range("E2").Value = '  +0.30%  '
$ws.Range("D3").Value = '''2.294.64'
$ws.Range("E3").Value = '  -0.15%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''97.34'
$ws.Range("E5").Value = '  +2.48%  '
$ws.Range("D6").Value = '''268.60'
$ws.Range("E6").Value = '  +0.19%  '
$ws.Range("E7").Value = '  -0.41%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = '''0.611'
$ws.Range("E9").Value = '  -1.55%  '
$ws.Range("D10").Value = '''45.62'
$ws.Range("E10").Value = '  +2.11%  '
$ws.Range("D11").Value = '''0.0936'
$ws.Range("E11").Value = '  +0.17%  '
$ws.Range("D12").Value = '''7.91'
$ws.Range("E12").Value = '  -2.18%  '
$ws.Range("D13").Value = '''0.107'
$ws.Range("E13").Value = '  +1.52%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = '''15.63'
$ws.Range("E14").Value = '  +2.38%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '''2.637.99'
$ws.Range("E15").Value = '  -0.18%  '
$ws.Range("D16").Value = '''0.856'
$ws.Range("E16").Value = '  +0.31%  '
$ws.Range("D17").Value = '''2.299.98'
$ws.Range("E17").Value = '  -0.40%  '
$ws.Range("D18").Value = '''43.736.00'
$ws.Range("E18").Value = '  +0.41%  '
$ws.Range("D19").Value = '''0.0000112'
$ws.Range("E19").Value = '  +4.05%  '
$ws.Range("D20").Value = '''6.18'
$ws.Range("E20").Value = '  -1.89%  '
$ws.Range("D21").Value = '''72.01'
$ws.Range("E21").Value = '  +1.19%  '
$ws.Range("D22").Value = '''2.53'
$ws.Range("E22").Value = '  +11.39%  '
$ws.Range("D23").Value = '''233.06'
$ws.Range("E23").Value = '  -1.76%  '
$ws.Range("D24").Value = '''9.06'
$ws.Range("E24").Value = '  -4.84%  '
$ws.Range("D25").Value = '''2.64'
$ws.Range("E25").Value = '  +6.49%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").Value = '''11.29'
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("E28").Value = '  +2.37%  '
$ws.Range("D29").Value = '''2.28'
$ws.Range("E29").Value = '  -1.17%  '
$ws.Range("D30").Value = '''39.09'
$ws.Range("E30").Value = '  +1.45%  '
$ws.Range("D31").Value = '''175.20'
$ws.Range("E31").Value = '  +2.09%  '
$ws.Range("D32").Value = '''21.93'
$ws.Range("E32").Value = '  -1.78%  '
$ws.Range("D33").Value = '''0.0902'
$ws.Range("E33").Value = '  +0.48%  '
$ws.Range("D34").Value = '''5.41'
$ws.Range("E34").Value = '  -2.15%  '
$ws.Range("E35").Value = '  -0.14%  '
$ws.Range("D36").Value = '''4.52'
$ws.Range("E36").Value = '  +1.79%  '
$ws.Range("E37").Value = '  -0.24%  '
$ws.Range("D38").Value = '''0.0353'
$ws.Range("E38").Value = '  -1.10%  '
$ws.Range("D39").Value = '''3.38'
$ws.Range("E39").Value = '  -2.46%  '
$ws.Range("E40").Value = '  +2.85%  '
$ws.Range("E41").Value = '  +1.87%  '
$ws.Range("D42").Value = '''12.25'
$ws.Range("E42").Value = '  +1.66%  '
$ws.Range("E43").Value = '  +0.73%  '
$ws.Range("D44").Value = '''64.55'
$ws.Range("E44").Value = '  +4.75%  '
$ws.Range("D45").Value = '''8.79'
$ws.Range("E45").Value = '  -2.43%  '
$ws.Range("B46").Value = 'THORChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D46").Value = '''5.16'
$ws.Range("E46").Value = '  -5.11%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '''0.102'
$ws.Range("E47").Value = '  -0.15%  '
$ws.Range("E48").Value = '  -2.47%  '
$ws.Range("D49").Value = '''1.20'
$ws.Range("E49").Value = '  -0.60%  '
$ws.Range("D50").Value = '''1.52'
$ws.Range("E50").Value = '  +12.05%  '
$ws.Range("D51").Value = '''2.518.02'
$ws.Range("E51").Value = '  -0.11%  '
